# Apply cryptocurrency price/volume(1h) updates to Sheet1
# (values written as text via a leading apostrophe so they remain
# literal strings, matching the original inline-string cell type)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'308.50"
$ws.Range("E2").Formula = "'0.11%"
$ws.Range("E3").Formula = "'3.21%"
$ws.Range("D4").Formula = "'5.111"
$ws.Range("E4").Formula = "'0.15%"
$ws.Range("D5").Formula = "'0.07625"
$ws.Range("E5").Formula = "'-1.08%"
$ws.Range("D6").Formula = "'4.248"
$ws.Range("E6").Formula = "'0.41%"
$ws.Range("D7").Formula = "'1.602"
$ws.Range("E7").Formula = "'-0.31%"
$ws.Range("E8").Formula = "'2.08%"
$ws.Range("D9").Formula = "'0.9015"
$ws.Range("E9").Formula = "'0.66%"
$ws.Range("E10").Formula = "'11.54%"
$ws.Range("D11").Formula = "'0.1798"
$ws.Range("E11").Formula = "'3.94%"
$ws.Range("D12").Formula = "'0.09175"
$ws.Range("E12").Formula = "'1.65%"
$ws.Range("D13").Formula = "'0.04169"
$ws.Range("E13").Formula = "'-5.73%"
$ws.Range("E14").Formula = "'-0.13%"
$ws.Range("D15").Formula = "'0.001253"
$ws.Range("E15").Formula = "'-0.29%"
$ws.Range("D16").Formula = "'0.005808"
$ws.Range("E16").Formula = "'-1.90%"
$ws.Range("D17").Formula = "'3.339"
$ws.Range("E17").Formula = "'-0.41%"
$ws.Range("D18").Formula = "'0.3314"
$ws.Range("E18").Formula = "'-0.10%"
$ws.Range("D19").Formula = "'6.655"
$ws.Range("E19").Formula = "'-5.82%"
$ws.Range("D20").Formula = "'0.1365"
$ws.Range("E20").Formula = "'1.20%"
$ws.Range("D21").Formula = "'0.2705"
$ws.Range("E21").Formula = "'-2.12%"
$ws.Range("E22").Formula = "'-1.66%"
$ws.Range("E23").Formula = "'3.12%"
$ws.Range("D24").Formula = "'0.004104"
$ws.Range("E24").Formula = "'1.10%"
$ws.Range("E25").Formula = "'-0.12%"
$ws.Range("D38").Formula = "'0.02399"
$ws.Range("E38").Formula = "'2.46%"
$ws.Range("D39").Formula = "'0.05189"
$ws.Range("E39").Formula = "'0.13%"
$ws.Range("D40").Formula = "'0.007780"
$ws.Range("E40").Formula = "'-2.12%"
$ws.Range("D41").Formula = "'0.1299"
$ws.Range("E41").Formula = "'-1.43%"
$ws.Range("D42").Formula = "'0.007049"
$ws.Range("E42").Formula = "'11.95%"
$ws.Range("E43").Formula = "'-0.09%"
$ws.Range("D44").Formula = "'0.007734"
$ws.Range("E44").Formula = "'-6.03%"
$ws.Range("D45").Formula = "'0.3078"
$ws.Range("E45").Formula = "'-7.67%"
$ws.Range("D46").Formula = "'0.00006956"
$ws.Range("E46").Formula = "'6.80%"
$ws.Range("D48").Formula = "'0.05461"
$ws.Range("E48").Formula = "'1,444.35%"
